$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("npm3d")
$ws2.Rows.Item(16).Select()
